$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# MOSIP-21520: strip the stray leading space from the Kannada and Tamil
# hierarchy_level_name values (kan: " ಪ್ರದೇಶ"/" ನಗರ"/" ಕೋಡ್", tam: " நாடு"/" நகரம்")
$ws.Range("C21").Value = "ಪ್ರದೇಶ"
$ws.Range("C23").Value = "ನಗರ"
$ws.Range("C25").Value = "ಕೋಡ್"
$ws.Range("C32").Value = "நாடு"
$ws.Range("C35").Value = "நகரம்"

# Update the sheet's saved view/scroll position and selection
$win = $excel.ActiveWindow
$win.ScrollRow = 18
$win.ScrollColumn = 1
$ws.Range("C30").Select()
